# Add a new row of data (row 9, skipping row 8) to Sheet1 and move the
# active selection down past it, mirroring the authored diff:
#   - dimension grows from A1:C7 to A1:C9
#   - new row r="9" with cell A9 = 8
#   - selection moves from A8 (A8:XFD13) to A10 (A10)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data value in A9 (row 8 intentionally left empty, matching the diff).
$ws.Range("A9").Value = 8

# Update the current selection/active cell to A10, as in the target sheetView.
$ws.Range("A10").Select() | Out-Null
